$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 12) mirroring the existing table's shape
$ws.Range("A12").Value = 10001
$ws.Range("B12").Value = "ekijkhk"
$ws.Range("C12").Value = 5465
$ws.Range("D12").Value = "PRJ-02"

# Move the active selection as recorded after the edit
$ws.Range("D11").Select()
